$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F - 想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3645
$ws1.Range("F5").Value = 2233
$ws1.Range("F6").Value = 435
$ws1.Range("F9").Value = 93
$ws1.Range("F11").Value = 1345
$ws1.Range("F13").Value = 2004
$ws1.Range("F14").Value = 145

# Sheet "全部类型" updates (column F - 想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3645
$ws4.Range("F5").Value = 2233
$ws4.Range("F6").Value = 435
$ws4.Range("F10").Value = 93
$ws4.Range("F14").Value = 1345
$ws4.Range("F16").Value = 2004
$ws4.Range("F17").Value = 145
